# Refresh cryptos list snapshot (prices + 1h volume %) per upstream scrape.
# Row 13/14 and 40/41 swap rank order (ShibaInu<->Avalanche, Cosmos<->Kaspa).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep its literal text (Excel would otherwise
    # silently coerce strings like "1.00" or "0.0000229" into numbers),
    # then restore the default style so no stray formatting is left behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "66.845.42"
$ws.Range("E2").Value = "  -1.94%  "

Set-TextValue $ws.Range("D3") "3.208.57"
$ws.Range("E3").Value = "  -3.96%  "

Set-TextValue $ws.Range("D4") "0.998"
$ws.Range("E4").Value = "  +0.04%  "

Set-TextValue $ws.Range("D5") "579.73"
$ws.Range("E5").Value = "  -4.15%  "

Set-TextValue $ws.Range("D6") "138.69"
$ws.Range("E6").Value = "  -15.27%  "

$ws.Range("E7").Value = "  +0.15%  "

Set-TextValue $ws.Range("D8") "3.197.05"
$ws.Range("E8").Value = "  -4.15%  "

Set-TextValue $ws.Range("D9") "0.519"
$ws.Range("E9").Value = "  -10.39%  "

Set-TextValue $ws.Range("D10") "0.160"
$ws.Range("E10").Value = "  -14.05%  "

Set-TextValue $ws.Range("D11") "6.43"
$ws.Range("E11").Value = "  -4.27%  "

Set-TextValue $ws.Range("D12") "0.473"
$ws.Range("E12").Value = "  -11.15%  "

$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D13") "0.0000229"
$ws.Range("E13").Value = "  -11.25%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D14") "35.66"
$ws.Range("E14").Value = "  -14.67%  "

Set-TextValue $ws.Range("D15") "3.725.85"
$ws.Range("E15").Value = "  -3.86%  "

Set-TextValue $ws.Range("D16") "66.746.87"
$ws.Range("E16").Value = "  -2.02%  "

Set-TextValue $ws.Range("D17") "3.207.06"
$ws.Range("E17").Value = "  -3.55%  "

$ws.Range("E18").Value = "  -5.17%  "

Set-TextValue $ws.Range("D19") "6.70"
$ws.Range("E19").Value = "  -13.42%  "

Set-TextValue $ws.Range("D20") "492.79"
$ws.Range("E20").Value = "  -12.23%  "

Set-TextValue $ws.Range("D21") "14.16"
$ws.Range("E21").Value = "  -12.56%  "

$ws.Range("E22").Value = "  -11.98%  "

Set-TextValue $ws.Range("D23") "7.24"
$ws.Range("E23").Value = "  -14.83%  "

Set-TextValue $ws.Range("D24") "81.61"
$ws.Range("E24").Value = "  -9.72%  "

Set-TextValue $ws.Range("D25") "12.65"
$ws.Range("E25").Value = "  -11.47%  "

Set-TextValue $ws.Range("D26") "1.00"
$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("E27").Value = "  -12.02%  "

Set-TextValue $ws.Range("D28") "27.51"
$ws.Range("E28").Value = "  -11.94%  "

Set-TextValue $ws.Range("D29") "2.01"
$ws.Range("E29").Value = "  -12.05%  "

Set-TextValue $ws.Range("D30") "7.43"
$ws.Range("E30").Value = "  -9.33%  "

$ws.Range("E31").Value = "  -4.06%  "

Set-TextValue $ws.Range("D32") "2.46"
$ws.Range("E32").Value = "  -7.16%  "

Set-TextValue $ws.Range("D33") "1.00"
$ws.Range("E33").Value = "  -0.18%  "

Set-TextValue $ws.Range("D34") "54.37"
$ws.Range("E34").Value = "  -1.13%  "

Set-TextValue $ws.Range("D35") "6.00"
$ws.Range("E35").Value = "  -18.16%  "

Set-TextValue $ws.Range("D36") "487.32"
$ws.Range("E36").Value = "  -16.07%  "

$ws.Range("E37").Value = "  -16.39%  "

Set-TextValue $ws.Range("D38") "0.0412"
$ws.Range("E38").Value = "  -8.89%  "

Set-TextValue $ws.Range("D39") "0.0805"
$ws.Range("E39").Value = "  -11.38%  "

$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D40") "8.38"
$ws.Range("E40").Value = "  -15.69%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D41") "0.118"
$ws.Range("E41").Value = "  -13.62%  "

Set-TextValue $ws.Range("D42") "2.806.14"
$ws.Range("E42").Value = "  -8.96%  "

Set-TextValue $ws.Range("D43") "2.46"
$ws.Range("E43").Value = "  -16.78%  "

$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("E45").Value = "  -11.10%  "

Set-TextValue $ws.Range("D46") "121.12"
$ws.Range("E46").Value = "  -6.70%  "

Set-TextValue $ws.Range("D47") "24.68"
$ws.Range("E47").Value = "  -14.78%  "

Set-TextValue $ws.Range("D48") "2.01"
$ws.Range("E48").Value = "  -10.90%  "

Set-TextValue $ws.Range("D49") "0.0₃0524"
$ws.Range("E49").Value = "  -16.57%  "

$ws.Range("E50").Value = "  -10.49%  "

$ws.Range("E51").Value = "  -21.02%  "
